$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 76.1785659790039
$ws.Range("C2").Value = 7.902298927307129
$ws.Range("D2").Value = 40.212406158447266
$ws.Range("H2").Value = 6.25
